$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.328788876533508
$ws.Range("B1").Value = 2.462164640426636
$ws.Range("C1").Value = 5.968862533569336
$ws.Range("D1").Value = 1.841139554977417
$ws.Range("E1").Value = 1.275335550308228
